$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.667.58'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '3.087.49'
$ws.Range('E3').Value = '  +1.28%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').Value = '''515.34'
$ws.Range('E5').Value = '  -0.89%  '
$ws.Range('D6').Value = '''143.06'
$ws.Range('E6').Value = '  +2.39%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').Value = '''0.436'
$ws.Range('E8').Value = '  +0.20%  '
$ws.Range('D9').Value = '''7.33'
$ws.Range('E9').Value = '  -0.18%  '
$ws.Range('E10').Value = '  -0.95%  '
$ws.Range('D11').Value = '''0.374'
$ws.Range('E11').Value = '  +0.65%  '
$ws.Range('D12').Value = '3.618.01'
$ws.Range('E12').Value = '  +2.25%  '
$ws.Range('E13').Value = '  +2.16%  '
$ws.Range('D14').Value = '''25.82'
$ws.Range('E14').Value = '  -3.99%  '
$ws.Range('D15').Value = '''0.0000165'
$ws.Range('E15').Value = '  -2.56%  '
$ws.Range('D16').Value = '57.764.95'
$ws.Range('E16').Value = '  +0.35%  '
$ws.Range('D17').Value = '3.095.23'
$ws.Range('E17').Value = '  +2.61%  '
$ws.Range('D18').Value = '''6.13'
$ws.Range('E18').Value = '  -1.68%  '
$ws.Range('D19').Value = '''13.10'
$ws.Range('E19').Value = '  -2.64%  '
$ws.Range('D20').Value = '''8.22'
$ws.Range('E20').Value = '  +0.45%  '
$ws.Range('D21').Value = '''336.51'
$ws.Range('E21').Value = '  +2.10%  '
$ws.Range('E22').Value = '  -0.01%  '
$ws.Range('D23').Value = '''0.502'
$ws.Range('E23').Value = '  -1.39%  '
$ws.Range('D24').Value = '''65.60'
$ws.Range('E24').Value = '  +1.41%  '
$ws.Range('E25').Value = '  +4.61%  '
$ws.Range('E26').Value = '  +0.45%  '
$ws.Range('D27').Value = '0.0₃0932'
$ws.Range('E27').Value = '  +3.53%  '
$ws.Range('D28').Value = '''6.48'
$ws.Range('E28').Value = '  -3.10%  '
$ws.Range('D29').Value = '''7.10'
$ws.Range('E29').Value = '  -2.31%  '
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').Value = '''20.90'
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('E32').Value = '  -3.50%  '
$ws.Range('D33').Value = '''154.41'
$ws.Range('E33').Value = '  +0.20%  '
$ws.Range('D34').Value = '''28.47'
$ws.Range('E34').Value = '  +11.46%  '
$ws.Range('D35').Value = '''4.53'
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').Value = '''5.92'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('D37').Value = '''1.24'
$ws.Range('E37').Value = '  -2.27%  '
$ws.Range('D38').Value = '''0.0687'
$ws.Range('E38').Value = '  +2.02%  '
$ws.Range('D39').Value = '3.136.61'
$ws.Range('E39').Value = '  +1.96%  '
$ws.Range('D40').Value = '''36.94'
$ws.Range('E40').Value = '  -0.99%  '
$ws.Range('B41').Value = 'Mantle'
$ws.Range('C41').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D41').Value = '''0.674'
$ws.Range('E41').Value = '  +1.41%  '
$ws.Range('B42').Value = 'Filecoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D42').Value = '''3.87'
$ws.Range('E42').Value = '  -0.64%  '
$ws.Range('E43').Value = '  +0.16%  '
$ws.Range('D44').Value = '2.285.66'
$ws.Range('E44').Value = '  +4.31%  '
$ws.Range('D45').Value = '''0.0254'
$ws.Range('E45').Value = '  +3.96%  '
$ws.Range('E46').Value = '  -0.77%  '
$ws.Range('D47').Value = '''20.37'
$ws.Range('E47').Value = '  +2.93%  '
$ws.Range('D48').Value = '''0.949'
$ws.Range('E48').Value = '  -0.61%  '
$ws.Range('D49').Value = '''5.89'
$ws.Range('E49').Value = '  -4.65%  '
$ws.Range('D50').Value = '''0.0878'
$ws.Range('E50').Value = '  +1.65%  '
$ws.Range('D51').Value = '''0.692'
$ws.Range('E51').Value = '  +2.63%  '
